$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lines")

# Insert a new column D ("Date Refused") before the existing "Lease Reference"
# column, shifting everything from D..Q to E..R.
$ws.Columns.Item(4).Insert()

# New header cell.
$ws.Cells.Item(1, 4).Value = "Date Refused"

# Row 2: lease amendment state PROPOSED -> REFUSED, and record the refusal
# date (same as the "Date Signed" value already on the row).
$ws.Cells.Item(2, 2).Value = "REFUSED"
$ws.Cells.Item(2, 4).Value = $ws.Cells.Item(2, 3).Value()

# Row 4: lease amendment state SIGNED -> REFUSED, with a matching refusal date.
$ws.Cells.Item(4, 2).Value = "REFUSED"
$ws.Cells.Item(4, 4).Value = $ws.Cells.Item(4, 3).Value()

# Row 5: lease amendment state SIGNED -> REFUSED, with a matching refusal date.
$ws.Cells.Item(5, 2).Value = "REFUSED"
$ws.Cells.Item(5, 4).Value = $ws.Cells.Item(5, 3).Value()

# Row 6: lease amendment state SIGNED -> REFUSED, with a matching refusal date.
$ws.Cells.Item(6, 2).Value = "REFUSED"
$ws.Cells.Item(6, 4).Value = $ws.Cells.Item(6, 3).Value()

$ws.Range("D1").Select()
